# Generate Report for Handoff
# -----------------------------------------------------------------------
# The localization pipeline generated a new handoff round: a new source
# guid (f9f4158f-d243-4969-805b-10a63cf6cfd2) replaces the previous one
# (d9b1cfe7-82d9-4715-bc18-f66c5ccdcd8d) everywhere it is referenced
# (file names, xliff names, hyperlink text) and the handoff timestamps
# advance to reflect the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "d9b1cfe7-82d9-4715-bc18-f66c5ccdcd8d"
$newGuid = "f9f4158f-d243-4969-805b-10a63cf6cfd2"

$oldXliffHash = "bbf312484054020c5c8da06479a87bd91ddeaed9"
$newXliffHash = "3da9f2eb858093e1da52cf786f369d75579db3bb"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22249eb45bd76d3c1ef89fca05c977afbeb70fbd/e2e/"

# ---------------------------------------------------------------
# Overview sheet: File Name / Path And Name / Latest HO Xliff Generate Date
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-10-27 09:50:13"

# Rebuild the B2 hyperlink so the visible display text matches the new
# file name (the stored target relationship is refreshed too).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$repoBase$newGuid.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null
$wsOverview.Range("B2").Font.Underline = $true
$wsOverview.Range("B2").Font.Color = 15570276

# ---------------------------------------------------------------
# zh-cn sheet: Source File Name / Latest Handoff File / Latest Handoff Datetime
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newXliffHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-10-27 09:50:00"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "$repoBase$newGuid.md", [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
$wsZhCn.Range("A2").Font.Underline = $true
$wsZhCn.Range("A2").Font.Color = 15570276

# ---------------------------------------------------------------
# de-de sheet: Source File Name / Latest Handoff File / Latest Handoff Datetime
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newXliffHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-10-27 09:50:13"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "$repoBase$newGuid.md", [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
$wsDeDe.Range("A2").Font.Underline = $true
$wsDeDe.Range("A2").Font.Color = 15570276
